$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
("111","222","2","OUT","2024-11-06-17-50-06"),
("64","157","4","OUT","2024-11-06-17-50-34"),
("105","215","4","OUT","2024-11-06-17-50-35"),
("67","154","4","OUT","2024-11-06-17-50-40"),
("67","154","4","OUT","2024-11-06-17-50-40"),
("93","207","4","OUT","2024-11-06-17-50-40"),
("136","201","4","OUT","2024-11-06-17-50-41"),
("81","166","3","OUT","2024-11-06-17-51-05"),
("87","184","3","OUT","2024-11-06-17-51-09"),
("113","212","3","OUT","2024-11-06-17-51-10"),
("73","146","2","OUT","2024-11-06-17-51-35"),
("76","159","2","OUT","2024-11-06-17-51-39"),
("167","211","2","OUT","2024-11-06-17-51-39"),
("82","192","3","OUT","2024-11-06-17-52-51"),
("82","192","3","OUT","2024-11-06-17-52-51"),
("107","225","3","OUT","2024-11-06-17-52-52"),
("59","116","2","OUT","2024-11-06-17-53-23"),
("90","187","2","OUT","2024-11-06-17-53-23"),
("82","186","2","OUT","2024-11-06-17-54-53"),
("111","219","2","OUT","2024-11-06-17-59-40"),
("180","201","2","OUT","2024-11-06-17-59-42"),
("193","145","2","OUT","2024-11-06-17-59-42"),
("196","127","2","OUT","2024-11-06-17-59-42"),
("198","124","2","OUT","2024-11-06-17-59-42"),
("196","126","2","OUT","2024-11-06-17-59-43"),
("196","126","2","OUT","2024-11-06-17-59-43"),
("196","127","2","OUT","2024-11-06-17-59-44"),
("196","127","2","OUT","2024-11-06-17-59-44"),
("196","127","2","OUT","2024-11-06-17-59-44"),
("196","127","2","OUT","2024-11-06-17-59-44"),
("196","128","2","OUT","2024-11-06-17-59-45"),
("196","128","2","OUT","2024-11-06-17-59-46"),
("196","126","2","OUT","2024-11-06-18-00-19"),
("98","232","2","IN","2024-11-06-18-51-55"),
("62","121","4","OUT","2024-11-06-18-53-08"),
("73","377","4","OUT","2024-11-06-18-53-15"),
("138","265","2","IN","2024-11-06-18-54-00"),
("78","173","2","IN","2024-11-06-18-54-01"),
("73","153","2","OUT","2024-11-06-18-57-12"),
("141","228","2","IN","2024-11-06-19-16-43"),
("92","199","2","IN","2024-11-06-19-16-45"),
("90","197","2","IN","2024-11-06-19-16-46"),
("83","183","2","IN","2024-11-06-19-16-46"),
("83","180","2","IN","2024-11-06-19-16-46"),
("75","173","2","IN","2024-11-06-19-16-46"),
("81","171","2","IN","2024-11-06-19-16-47"),
("83","170","2","IN","2024-11-06-19-16-47"),
("84","211","2","IN","2024-11-06-19-17-45"),
("109","249","2","IN","2024-11-06-19-17-45"),
("77","217","3","IN","2024-11-06-19-19-23"),
("80","171","2","OUT","2024-11-06-19-20-28"),
("223","351","2","OUT","2024-11-06-19-20-29"),
("70","169","4","OUT","2024-11-06-19-20-43"),
("71","172","4","OUT","2024-11-06-19-20-44"),
("95","211","4","OUT","2024-11-06-19-20-44"),
("120","254","4","OUT","2024-11-06-19-20-45"),
("150","282","4","OUT","2024-11-06-19-20-45"),
("77","176","3","OUT","2024-11-06-19-21-21"),
("87","199","3","OUT","2024-11-06-19-21-22"),
("114","233","3","OUT","2024-11-06-19-21-22"),
("146","252","3","OUT","2024-11-06-19-21-23"),
("65","133","2","OUT","2024-11-06-19-21-53"),
("65","133","2","OUT","2024-11-06-19-21-53"),
("75","160","2","OUT","2024-11-06-19-21-54"),
("80","178","2","OUT","2024-11-06-19-21-54"),
("87","199","2","OUT","2024-11-06-19-21-54"),
("98","217","2","OUT","2024-11-06-19-21-55"),
("102","224","2","OUT","2024-11-06-19-21-55"),
("114","241","2","OUT","2024-11-06-19-21-55"),
("118","247","2","OUT","2024-11-06-19-21-55"),
("192","273","2","OUT","2024-11-06-19-21-56"),
("100","245","2","IN","2024-11-06-19-22-53")
)

$startRow = 59
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = [double]$vals[0]
    $ws.Cells.Item($row, 2).Value = [double]$vals[1]
    $ws.Cells.Item($row, 3).Value = [double]$vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
}
